# Correção do artefato 19 - reescreve as regras de negócio (RN-0001..RN-0005)
# como um novo conjunto RN-0001..RN-0006, com texto revisado e alinhamento
# justificado em cada parágrafo.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParagraphXml($paraIndex, $bodyXml) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $full = $pkgHeader + '<w:body>' + $bodyXml + '</w:body>' + $pkgFooter
    $r.InsertXML($full) | Out-Null
}

# Parágrafo 2 (antes vazio) -> RN-0001
$p2 = '<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="000000"/><w:lang w:eastAsia="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RN-0001:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:color w:val="000000"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/><w:lang w:eastAsia="pt-BR"/></w:rPr><w:t>O parceiro deve responder em até 7 dias a partir do primeiro contato.</w:t></w:r></w:p>'

# Parágrafo 3 (era RN-0001) -> RN-0002
$p3 = '<w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RN-0002:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">O cliente pode solicitar </w:t></w:r><w:r><w:t xml:space="preserve">adesão da excursão </w:t></w:r><w:r><w:t>em até um dia antes da data de embarque.</w:t></w:r></w:p>'

# Parágrafo 4 (era RN-0002) -> RN-0003
$p4 = '<w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RN-</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>0</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">003: </w:t></w:r><w:r><w:t>O</w:t></w:r><w:r><w:t xml:space="preserve"> cliente</w:t></w:r><w:r><w:t xml:space="preserve"> deve </w:t></w:r><w:r><w:t>realizar pagamento</w:t></w:r><w:r><w:t xml:space="preserve"> em até 48hrs úteis após a compra do pacote</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>'

# Parágrafo 5 (era RN-0003) -> RN-0004
$p5 = '<w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RN-0</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>0</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>0</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>O cliente pode solicitar cancelamento da excursão com até 3 dias antes da data de embarque.</w:t></w:r></w:p>'

# Parágrafo 6 (era RN-0004) -> RN-0005
$p6 = '<w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RN-0</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>0</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">5: </w:t></w:r><w:r><w:t xml:space="preserve">O cliente deve </w:t></w:r><w:r><w:t>apresentar o recibo de pagamento</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>para que a vaga seja garantida.</w:t></w:r></w:p>'

# Parágrafo 7 (era RN-0005) -> RN-0006
$p7 = '<w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>RN-0</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>0</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>6</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">O gerente realiza estorno </w:t></w:r><w:r><w:t>de 100% do valor</w:t></w:r><w:r><w:t xml:space="preserve"> em até dois dias úteis após a solicitação do cancelamento. Caso a solicitação seja feita fora do prazo será estornado apenas 70% do valor pago.</w:t></w:r></w:p>'

Set-ParagraphXml 2 $p2
Set-ParagraphXml 3 $p3
Set-ParagraphXml 4 $p4
Set-ParagraphXml 5 $p5
Set-ParagraphXml 6 $p6
Set-ParagraphXml 7 $p7

Write-Host "Regras de negocio atualizadas."
